$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells hold numeric-looking text (e.g. "42.20", "5.141") that must
# stay text (matching the source data) instead of being auto-converted to numbers,
# which would silently drop significant trailing zeros / change the stored type.
$textCells = @("D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D26", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '21.618.77'
$ws.Range('D3').Value = '1.532.46'
$ws.Range('E3').Value = '  -1.46%  '
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('D6').Value = '288.77'
$ws.Range('E6').Value = '  -0.01%  '
$ws.Range('D7').Value = '0.3955'
$ws.Range('E7').Value = '  +0.73%  '
$ws.Range('D8').Value = '0.3162'
$ws.Range('E8').Value = '  -1.83%  '
$ws.Range('D9').Value = '42.20'
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('D10').Value = '0.07156'
$ws.Range('E10').Value = '  -2.54%  '
$ws.Range('D11').Value = '1.051'
$ws.Range('E11').Value = '  -5.91%  '
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').Value = '5.701'
$ws.Range('E13').Value = '  +0.24%  '
$ws.Range('D14').Value = '18.31'
$ws.Range('E14').Value = '  -3.93%  '
$ws.Range('D15').Value = '6.593'
$ws.Range('E15').Value = '  -2.60%  '
$ws.Range('D16').Value = '1.542.72'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').Value = '0.00001087'
$ws.Range('E17').Value = '  -3.10%  '
$ws.Range('D18').Value = '0.06592'
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('D19').Value = '83.60'
$ws.Range('E19').Value = '  -1.72%  '
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').Value = '6.116'
$ws.Range('E21').Value = '  -4.21%  '
$ws.Range('E22').Value = '  -2.96%  '
$ws.Range('D23').Value = '10.62'
$ws.Range('E23').Value = '  -7.18%  '
$ws.Range('D24').Value = '2.361'
$ws.Range('E24').Value = '  +0.81%  '
$ws.Range('D25').Value = '21.614.92'
$ws.Range('E25').Value = '  -1.98%  '
$ws.Range('D26').Value = '2.330'
$ws.Range('E26').Value = '  -7.62%  '
$ws.Range('D27').Value = '149.16'
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('D28').Value = '18.32'
$ws.Range('E28').Value = '  -2.73%  '
$ws.Range('D29').Value = '4.835'
$ws.Range('E29').Value = '  -0.49%  '
$ws.Range('D30').Value = '1.717.25'
$ws.Range('E30').Value = '  -0.93%  '
$ws.Range('D31').Value = '116.79'
$ws.Range('E31').Value = '  -3.41%  '
$ws.Range('D32').Value = '6.050'
$ws.Range('E32').Value = '  +4.58%  '
$ws.Range('D33').Value = '0.9278'
$ws.Range('E33').Value = '  -13.64%  '
$ws.Range('D34').Value = '0.08108'
$ws.Range('E34').Value = '  -1.59%  '
$ws.Range('B35').Value = 'InternetComputer(DFINITY)'
$ws.Range('C35').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D35').Value = '5.141'
$ws.Range('E35').Value = '  -1.01%  '
$ws.Range('B36').Value = 'FraxShare'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D36').Value = '8.404'
$ws.Range('E36').Value = '  -10.28%  '
$ws.Range('D37').Value = '0.05970'
$ws.Range('E37').Value = '  -4.23%  '
$ws.Range('D38').Value = '0.02209'
$ws.Range('E38').Value = '  -3.53%  '
$ws.Range('D39').Value = '1.450'
$ws.Range('E39').Value = '  -13.32%  '
$ws.Range('D40').Value = '0.2025'
$ws.Range('E40').Value = '  -3.96%  '
$ws.Range('E41').Value = '  -2.99%  '
$ws.Range('D42').Value = '10.93'
$ws.Range('E42').Value = '  +1.58%  '
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').Value = '0.5765'
$ws.Range('E44').Value = '  -2.62%  '
$ws.Range('D45').Value = '12.92'
$ws.Range('E45').Value = '  -3.69%  '
$ws.Range('D46').Value = '3.710'
$ws.Range('E46').Value = '  -0.49%  '
$ws.Range('D47').Value = '0.5490'
$ws.Range('E47').Value = '  -3.80%  '
$ws.Range('D48').Value = '1.164'
$ws.Range('E48').Value = '  +1.31%  '
$ws.Range('D49').Value = '1.868'
$ws.Range('E49').Value = '  -2.96%  '
$ws.Range('D50').Value = '115.77'
$ws.Range('E50').Value = '  -2.11%  '
$ws.Range('D51').Value = '0.06677'
$ws.Range('E51').Value = '  -2.88%  '
